$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '47.271.32'
$ws.Range("E2").Value = '  +0.11%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.489.90'
$ws.Range("E3").Value = '  +0.12%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.19'
$ws.Range("E5").Value = '  -0.45%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.34'
$ws.Range("E6").Value = '  +3.41%  '

# Row 7
$ws.Range("E7").Value = '  +0.57%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.540'
$ws.Range("E9").Value = '  +0.56%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.21'
$ws.Range("E10").Value = '  +5.02%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0812'
$ws.Range("E11").Value = '  +0.16%  '

# Row 12
$ws.Range("E12").Value = '  +0.69%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.37'
$ws.Range("E13").Value = '  +0.79%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.16'
$ws.Range("E14").Value = '  -0.09%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.878.50'
$ws.Range("E15").Value = '  +0.20%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.491.50'
$ws.Range("E16").Value = '  +0.84%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.844'
$ws.Range("E17").Value = '  +0.58%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.164.90'
$ws.Range("E18").Value = '  +0.12%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.29'
$ws.Range("E19").Value = '  +5.13%  '

# Row 20
$ws.Range("E20").Value = '  +1.49%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0942'
$ws.Range("E21").Value = '  +1.08%  '

# Row 22
$ws.Range("E22").Value = '  +15.11%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.44'
$ws.Range("E23").Value = '  -0.03%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '245.89'
$ws.Range("E24").Value = '  -1.53%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.55'
$ws.Range("E25").Value = '  -0.29%  '

# Row 26
$ws.Range("E26").Value = '  +0.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.72'
$ws.Range("E27").Value = '  -1.34%  '

# Row 28
$ws.Range("E28").Value = '  +3.68%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.98'
$ws.Range("E29").Value = '  -0.86%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.139'
$ws.Range("E30").Value = '  +5.00%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.63'
$ws.Range("E31").Value = '  -0.64%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.71'
$ws.Range("E32").Value = '  +0.62%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.19'
$ws.Range("E33").Value = '  +1.37%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.33'
$ws.Range("E34").Value = '  -0.84%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0783'
$ws.Range("E35").Value = '  +0.29%  '

# Row 36
$ws.Range("E36").Value = '  +0.14%  '

# Row 37
$ws.Range("E37").Value = '  +2.38%  '

# Row 38
$ws.Range("E38").Value = '  +1.46%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.94'
$ws.Range("E39").Value = '  -1.32%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.01'
$ws.Range("E40").Value = '  +8.36%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.112'
$ws.Range("E41").Value = '  +0.57%  '

# Row 42
$ws.Range("E42").Value = '  +0.00%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '118.24'
$ws.Range("E43").Value = '  -1.84%  '

# Row 44
$ws.Range("E44").Value = '  +0.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.998.07'
$ws.Range("E45").Value = '  +2.57%  '

# Row 46
$ws.Range("E46").Value = '  +1.82%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.01'
$ws.Range("E47").Value = '  -3.98%  '

# Row 48
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.80'
$ws.Range("E48").Value = '  +0.44%  '

# Row 49
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.14'
$ws.Range("E49").Value = '  -0.75%  '

# Row 50
$ws.Range("E50").Value = '  -4.39%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.48'
$ws.Range("E51").Value = '  +3.39%  '
